$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.228.36'
$ws.Range("D3").Value = '2.242.73'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.67'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.629'
$ws.Range("E6").Value = '  -2.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.45'
$ws.Range("E7").Value = '  -3.18%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  -4.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.13'
$ws.Range("E10").Value = '  +6.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0942'
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.16'
$ws.Range("E12").Value = '  -1.64%  '
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.51'
$ws.Range("E14").Value = '  -2.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.852'
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").Value = '2.252.57'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '42.092.49'
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("D18").Value = '0.0₃0984'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.16'
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.13'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '231.48'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.21'
$ws.Range("E22").Value = '  +4.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.81'
$ws.Range("E23").Value = '  +38.50%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.44'
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.63'
$ws.Range("E26").Value = '  -4.59%  '
$ws.Range("E27").Value = '  -1.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.16'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.17'
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0820'
$ws.Range("E31").Value = '  -3.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.30'
$ws.Range("E32").Value = '  +3.28%  '
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.125'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.19'
$ws.Range("E35").Value = '  +10.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.48'
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0314'
$ws.Range("E37").Value = '  +3.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '13.71'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.75'
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '62.24'
$ws.Range("E41").Value = '  +2.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.204'
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.50'
$ws.Range("E43").Value = '  -3.23%  '
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.66'
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.16'
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.23'
$ws.Range("E49").Value = '  -5.22%  '
$ws.Range("E50").Value = '  +2.49%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.70'
$ws.Range("E51").Value = '  +0.51%  '
